$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> (new Price text, new Volume(1h) text, ForceText flag)
# Empty string means "no change for this column" (only one of D/E changed in the diff)
# ForceText = $true means the new Price value looks like a plain number and must be
# kept as text (matching the original inlineStr cell type) by pre-formatting as Text.
$updates = @(
    [PSCustomObject]@{ Row = 2; D = '61.302.91'; E = '  -4.76%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 3; D = '3.312.86'; E = '  -5.30%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 4; D = '1.00'; E = '  +0.04%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 5; D = '567.55'; E = '  -3.76%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 6; D = '126.57'; E = '  -5.78%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 7; D = ''; E = '  -0.05%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 8; D = '3.312.97'; E = '  -5.26%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 9; D = '0.475'; E = '  -2.75%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 10; D = '7.17'; E = '  -5.77%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 11; D = ''; E = '  -5.95%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 12; D = '0.373'; E = '  -4.90%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 13; D = '3.878.23'; E = '  -5.25%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 14; D = '0.118'; E = '  -1.41%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 15; D = '3.319.10'; E = '  -5.14%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 16; D = '0.0000168'; E = '  -7.09%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 17; D = '24.69'; E = '  -3.68%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 18; D = '61.407.00'; E = '  -4.61%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 19; D = '9.09'; E = '  -9.50%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 20; D = '5.57'; E = '  -3.73%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 21; D = '13.10'; E = '  -3.46%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 22; D = '356.82'; E = '  -8.38%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 23; D = '0.551'; E = '  -5.42%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 24; D = ''; E = '  -0.04%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 25; D = '3.443.93'; E = '  -5.33%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 26; D = '70.36'; E = '  -5.34%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 27; D = ''; E = '  -7.68%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 28; D = ''; E = '  +0.42%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 29; D = '7.19'; E = '  -2.83%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 30; D = '1.46'; E = '  -1.06%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 31; D = '7.87'; E = '  -3.82%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 32; D = '2.10'; E = '  -7.42%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 34; D = '0.147'; E = '  -5.68%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 35; D = '3.337.29'; E = '  -5.37%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 36; D = '5.46'; E = '  +1.95%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 37; D = '22.27'; E = '  -4.75%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 38; D = '6.72'; E = '  -3.21%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 39; D = '162.38'; E = '  -1.85%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 40; D = '1.49'; E = '  -3.79%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 41; D = '0.0754'; E = '  -4.16%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 42; D = '1.00'; E = '  +0.14%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 43; D = '40.89'; E = '  -2.67%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 44; D = '0.747'; E = '  -7.64%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 45; D = ''; E = '  -4.95%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 46; D = ''; E = '  -6.15%  '; ForceText = $false }
    [PSCustomObject]@{ Row = 47; D = '1.52'; E = '  -7.97%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 48; D = '22.24'; E = '  -10.18%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 49; D = '6.59'; E = '  -3.44%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 50; D = '0.854'; E = '  -7.76%  '; ForceText = $true }
    [PSCustomObject]@{ Row = 51; D = '2.190.92'; E = '  -8.87%  '; ForceText = $false }
)

foreach ($u in $updates) {
    if ($u.D -ne '') {
        $cell = $ws.Cells.Item($u.Row, 4)
        if ($u.ForceText) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    if ($u.E -ne '') {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
